# Add a new "updated_at" column (G) that stamps each data row with the
# current date/time via a volatile NOW() formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in G1
$ws.Range("G1").Value = "updated_at"

# Give the new column a date-time display format, then fill it with
# =NOW() — first cell gets the formula directly, the rest share it
# (mirrors how Excel itself records a dragged-down formula as a shared
# formula group).
$ws.Range("G2:G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("G2").Formula = "=NOW()"
$ws.Range("G3:G5").Formula = "=NOW()"

# Size the new column to fit its contents.
$ws.Columns.Item(7).AutoFit()

# Keep the sheet's recorded selection in sync with where the edit left it.
$ws.Range("G8").Select() | Out-Null

# Touch the page setup so a print orientation is recorded for the sheet.
$ws.PageSetup.Orientation = 1
